$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 change
$ws.Range("K2").Value = 1.8

# Row 4 changes
$ws.Range("G4").Value = 2.5
$ws.Range("I4").Value = 2.7
$ws.Range("L4").Value = 3.6
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1.73
$ws.Range("W4").Value = 7
$ws.Range("Y4").Value = 11
$ws.Range("Z4").Value = 26
$ws.Range("AJ4").Value = 11
$ws.Range("AL4").Value = 26
$ws.Range("AN4").Value = 4.5
$ws.Range("AX4").Value = 17
$ws.Range("AY4").Value = 29
$ws.Range("AZ4").Value = 51
$ws.Range("BB4").Value = 251
